$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 11.93779999999999
$ws.Range("D3").Value = -5.996099999999999
$ws.Range("A4").Value = -21.19120000000002
$ws.Range("B4").Value = 4.674900000000005
$ws.Range("D4").Value = -7.638100000000001
$ws.Range("B5").Value = 5.5404
$ws.Range("E5").Value = 13.09339999999999
$ws.Range("A6").Value = -21.57190000000001
$ws.Range("B6").Value = 5.551899999999998
$ws.Range("A7").Value = -21.4543
$ws.Range("A8").Value = -21.36760000000002
$ws.Range("B8").Value = 5.1829
$ws.Range("D9").Value = -8.367899999999995
$ws.Range("D11").Value = -8.267499999999998
$ws.Range("D14").Value = -8.352799999999997
$ws.Range("A16").Value = -21.37740000000002
$ws.Range("B16").Value = 5.150000000000003
$ws.Range("D18").Value = -7.866799999999995
$ws.Range("A20").Value = -22.20780000000003
$ws.Range("E20").Value = 13.29619999999998
$ws.Range("A21").Value = -20.50699999999999
$ws.Range("B22").Value = 5.094000000000005
$ws.Range("D25").Value = -7.868399999999999
